$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.382.00'
$ws.Range('E2').Value = '  -3.32%  '
$ws.Range('D3').Value = '3.167.02'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'607.42"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = "'148.28"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.22%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.163.33'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('E9').Value = '  -3.68%  '
$ws.Range('D10').Value = "'0.152"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.96%  '
$ws.Range('D11').Value = "'5.52"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.41%  '
$ws.Range('D12').Value = "'0.478"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.60%  '
$ws.Range('D13').Value = "'0.0000253"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.48%  '
$ws.Range('D14').Value = "'35.75"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.00%  '
$ws.Range('D15').Value = '3.682.28'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = '64.359.16'
$ws.Range('E16').Value = '  -3.42%  '
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '3.161.59'
$ws.Range('E18').Value = '  -2.65%  '
$ws.Range('E19').Value = '  -5.51%  '
$ws.Range('D20').Value = "'482.61"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.99%  '
$ws.Range('D21').Value = "'14.75"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('E22').Value = '  -5.10%  '
$ws.Range('D23').Value = "'7.78"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').Value = "'13.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.27%  '
$ws.Range('D25').Value = "'83.81"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.23%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  -4.96%  '
$ws.Range('D28').Value = "'8.52"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('D29').Value = "'2.20"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.39%  '
$ws.Range('D30').Value = "'6.81"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  -17.37%  '
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = "'26.35"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.13%  '
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = "'55.03"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').Value = "'6.02"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.36%  '
$ws.Range('D38').Value = '0.0₃0732'
$ws.Range('E38').Value = '  -7.70%  '
$ws.Range('D39').Value = "'455.30"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.50%  '
$ws.Range('D40').Value = "'2.95"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.87%  '
$ws.Range('D41').Value = "'0.0398"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.00%  '
$ws.Range('E42').Value = '  -4.06%  '
$ws.Range('E43').Value = '  -8.28%  '
$ws.Range('D44').Value = '2.854.47'
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('E45').Value = '  -8.56%  '
$ws.Range('D46').Value = "'2.29"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.70%  '
$ws.Range('D47').Value = "'26.51"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.29%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = "'2.32"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.05%  '
$ws.Range('E50').Value = '  -4.57%  '
$ws.Range('D51').Value = "'119.77"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.44%  '
